$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.170.90'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").Value = '2.247.55'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.97'
$ws.Range("E5").Value = '  -2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.90'
$ws.Range("E7").Value = '  +3.75%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.628'
$ws.Range("E9").Value = '  -2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.18'
$ws.Range("E10").Value = '  +5.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("E11").Value = '  -3.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.11'
$ws.Range("E12").Value = '  -4.71%  '
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '2.584.24'
$ws.Range("E14").Value = '  -2.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.75'
$ws.Range("E15").Value = '  -4.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.859'
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '2.244.94'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '42.089.44'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.97'
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E21").Value = '  -3.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.28'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.87'
$ws.Range("E23").Value = '  -3.08%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.33'
$ws.Range("E25").Value = '  -3.18%  '
$ws.Range("E26").Value = '  -7.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.55'
$ws.Range("E27").Value = '  +18.58%  '
$ws.Range("E28").Value = '  -5.27%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.16'
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.97'
$ws.Range("E30").Value = '  +1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.58'
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0837'
$ws.Range("E32").Value = '  -0.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.69'
$ws.Range("E33").Value = '  +5.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.120'
$ws.Range("E34").Value = '  -5.50%  '
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.50'
$ws.Range("E36").Value = '  -2.66%  '
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.41'
$ws.Range("E38").Value = '  +3.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0302'
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  -7.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '112.70'
$ws.Range("E42").Value = '  +7.33%  '
$ws.Range("E43").Value = '  -7.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.95'
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.68'
$ws.Range("E45").Value = '  -5.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0992'
$ws.Range("E46").Value = '  -4.53%  '
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("E48").Value = '  -4.57%  '
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.28'
$ws.Range("E50").Value = '  -12.51%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.30'
$ws.Range("E51").Value = '  -0.68%  '
